$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Agrega parte 1 de nuevos estado de cuenta": the account statement now
# covers only the first worker/period; drop the rows that listed
# PAULA ANDREA CABARCAS VELEZ for periods 2506 and 2505. Deleting the
# entire rows also shifts the signature block up and refreshes the
# shared-string / used-range bookkeeping automatically.
$ws.Range("B17:B18").EntireRow.Delete()

# "Actualiza base de datos EC": refresh the summary counters to match
# the now-single worker / single overdue period left on the sheet.
$ws.Range("E11").Value = 37960
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
